$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Maven file path value in G4 from "c:\testupdate" to "c:\test"
$ws.Range("G4").Value = "c:\test"

# Update selection to G4 to match the saved view state
$ws.Range("G4").Select()
